$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-19 -> 2023-09-20, i.e. serial 45188 -> 45189) for every data
# row (rows 2 through 201).
$ws.Range("C2:C201").Value = 45189
